$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column P held the conference as text ("West"/"East"); turn it into a binary
# variable: West -> 0, East -> 1. Column Q's formula double-counted the
# ThreePM column (E) when it shouldn't have - fix it to just C - D.
for ($r = 2; $r -le 31; $r++) {
    $pCell = $ws.Cells.Item($r, 16)   # column P (Conf)
    $conf = $pCell.Value2

    if ($conf -eq "West") {
        $pCell.Value = 0
    } elseif ($conf -eq "East") {
        $pCell.Value = 1
    }

    $qCell = $ws.Cells.Item($r, 17)   # column Q
    $qCell.Formula = "=C$r-(D$r)"
}

# Restore the view state recorded in the workbook: scrolled so row 14 is at
# the top, with P2 as the active/selected cell.
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 14
$aw.ScrollColumn = 1
$ws.Range("P2").Select() | Out-Null
